# Insert a new weekly record row for "Rabanito" (Vega Central Mapocho de Santiago)
# at row 280, shifting the existing rows 280:302 down to 281:303.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 280 (pushes old row 280 .. 302 down to 281 .. 303)
$ws.Rows.Item(280).Insert()

# Populate the new row 280 with the new data record
$ws.Cells.Item(280, 1).Value = 9
$ws.Cells.Item(280, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(280, 3).Value = "Metropolitana"
$ws.Range("D280").Value = 44783
$ws.Cells.Item(280, 5).Value = 13
$ws.Cells.Item(280, 6).Value = 300000001
$ws.Cells.Item(280, 7).Value = "Rabanito"
$ws.Cells.Item(280, 8).Value = "Sin especificar"
$ws.Cells.Item(280, 9).Value = "Primera"
$ws.Cells.Item(280, 10).Value = 10400
$ws.Cells.Item(280, 11).Value = 2500
$ws.Cells.Item(280, 12).Value = 3000
$ws.Cells.Item(280, 13).Value = 2750
$ws.Cells.Item(280, 14).Value = "`$/cien unidades (volumen en unidades)"
$ws.Cells.Item(280, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(280, 16).Value = 28
$ws.Cells.Item(280, 17).Value = 100
$ws.Cells.Item(280, 18).Value = "Hortaliza"
